$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.430.16'
$ws.Range("E2").Value = '  -2.67%  '
$ws.Range("D3").Value = '1.784.37'
$ws.Range("E3").Value = '  -2.43%  '
$ws.Range("D5").Value = '229.84'
$ws.Range("E5").Value = '  -1.87%  '
$ws.Range("D6").Value = '0.5849'
$ws.Range("E6").Value = '  -2.53%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '0.2741'
$ws.Range("E8").Value = '  -0.66%  '
$ws.Range("D9").Value = '23.06'
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("D10").Value = '0.06654'
$ws.Range("E10").Value = '  -4.69%  '
$ws.Range("D11").Value = '0.07532'
$ws.Range("E11").Value = '  -1.12%  '
$ws.Range("D12").Value = '1.785.96'
$ws.Range("E12").Value = '  -2.51%  '
$ws.Range("D13").Value = '4.748'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("D14").Value = '0.6048'
$ws.Range("E14").Value = '  -3.65%  '
$ws.Range("D15").Value = '2.026.55'
$ws.Range("E15").Value = '  -2.41%  '
$ws.Range("D16").Value = '74.76'
$ws.Range("E16").Value = '  -4.60%  '
$ws.Range("D17").Value = '0.000008564'
$ws.Range("E17").Value = '  -11.44%  '
$ws.Range("D18").Value = '28.388.07'
$ws.Range("E18").Value = '  -1.55%  '
$ws.Range("D19").Value = '5.360'
$ws.Range("E19").Value = '  -6.12%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").Value = '205.90'
$ws.Range("E21").Value = '  -6.92%  '
$ws.Range("D22").Value = '11.33'
$ws.Range("E22").Value = '  -1.82%  '
$ws.Range("D23").Value = '6.731'
$ws.Range("E23").Value = '  -1.87%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = '151.83'
$ws.Range("E25").Value = '  -2.30%  '
$ws.Range("D26").Value = '8.076'
$ws.Range("E26").Value = '  +1.32%  '
$ws.Range("D27").Value = '0.1245'
$ws.Range("E27").Value = '  -3.67%  '
$ws.Range("D28").Value = '16.21'
$ws.Range("E28").Value = '  -2.06%  '
$ws.Range("D29").Value = '1.405'
$ws.Range("E29").Value = '  -3.29%  '
$ws.Range("D30").Value = '0.06079'
$ws.Range("E30").Value = '  -7.27%  '
$ws.Range("D31").Value = '1.413'
$ws.Range("E31").Value = '  -1.67%  '
$ws.Range("D32").Value = '3.750'
$ws.Range("E32").Value = '  -2.20%  '
$ws.Range("D33").Value = '3.758'
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("D34").Value = '1.664'
$ws.Range("E34").Value = '  -3.44%  '
$ws.Range("D35").Value = '1.037'
$ws.Range("E35").Value = '  -5.38%  '
$ws.Range("D36").Value = '0.6357'
$ws.Range("E36").Value = '  -1.43%  '
$ws.Range("D37").Value = '2.504'
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("D38").Value = '2.688'
$ws.Range("E38").Value = '  -1.68%  '
$ws.Range("D39").Value = '1.140.48'
$ws.Range("E39").Value = '  -2.87%  '
$ws.Range("D40").Value = '0.01666'
$ws.Range("E40").Value = '  -4.68%  '
$ws.Range("D41").Value = '6.258'
$ws.Range("E41").Value = '  -4.19%  '
$ws.Range("D42").Value = '0.8719'
$ws.Range("E42").Value = '  -2.85%  '
$ws.Range("D44").Value = '100.34'
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").Value = '1.935.40'
$ws.Range("E45").Value = '  -2.42%  '
$ws.Range("D46").Value = '59.68'
$ws.Range("E46").Value = '  -3.97%  '
$ws.Range("E47").Value = '  -2.86%  '
$ws.Range("D48").Value = '8.339'
$ws.Range("E48").Value = '  -1.51%  '
$ws.Range("D49").Value = '1.561'
$ws.Range("E49").Value = '  -1.92%  '
$ws.Range("D50").Value = '0.05414'
$ws.Range("E50").Value = '  -3.70%  '
$ws.Range("E51").Value = '  -1.82%  '
